# Update the "dSF" column (F) values for the davidson_tucker workbook.
# Repull / push all data + recompute mean -> several dSF values change
# to diverge from dS0 (column E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = 5
    "F3"  = -1
    "F5"  = -5
    "F6"  = 2
    "F7"  = -1
    "F8"  = 1
    "F9"  = -3
    "F10" = -1
    "F11" = -3
    "F12" = 1
    "F13" = -10
    "F14" = 1
    "F15" = -3
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
